$d = $word.ActiveDocument

# --- Locate the paragraph that holds the original run text ---
$p1 = $d.Paragraphs(1)

# Original text is "fdsafdasfdsfdsfsdfdsfdfdf" (25 chars), all in one run.
# Target: three runs -> "F" / "dsafdasfdsfdsfsdfdsfdfdf" / "_1111",
# plus a new empty paragraph right after this one.

$paraStart = $p1.Range.Start
$paraTextEnd = $p1.Range.End - 1   # exclude the paragraph mark

# Step 1: capitalize the first letter "f" -> "F" (plain text substitution).
$rFirst = $d.Range($paraStart, $paraStart + 1)
$rFirst.Text = "F"

# Step 2: append "_1111" right after the existing text, before the paragraph mark.
$p1 = $d.Paragraphs(1)
$paraTextEnd = $p1.Range.End - 1
$rTailInsert = $d.Range($paraTextEnd, $paraTextEnd)
$rTailInsert.InsertAfter("_1111")

# Step 3: insert a new, empty paragraph right after this paragraph.
$p1 = $d.Paragraphs(1)
$paraEnd = $p1.Range.End
$rNewPara = $d.Range($paraEnd, $paraEnd)
$rNewPara.InsertParagraphAfter()

# Step 4: split "F" into its own run (toggle a character property to force a run
# boundary, then restore it so the visible formatting is unchanged).
$rSplit1 = $d.Range($paraStart, $paraStart + 1)
$rSplit1.Font.Bold = 1
$rSplit1.Font.Bold = 0

# Step 5: split "_1111" into its own trailing run the same way.
$p1 = $d.Paragraphs(1)
$paraTextEnd = $p1.Range.End - 1
$rSplit2 = $d.Range($paraTextEnd - 5, $paraTextEnd)
$rSplit2.Font.Bold = 1
$rSplit2.Font.Bold = 0

Write-Output $d.Content.Text
Write-Output $d.Paragraphs.Count
